# Rename the worksheet tab to reflect the updated dataset naming
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "1_Mohgaon_S_RP_RL_data_univaria"

# Build the replacement data block (rows 2-100, columns A-E)
$data = New-Object 'object[,]' 99,5
$data[0,0] = 2357.6962121279598
$data[0,1] = 1792.76996480348
$data[0,2] = 2934.43789662779
$data[0,3] = 1617.5642688486
$data[0,4] = 3118.24203173572
$data[1,0] = 3196.1695425937701
$data[1,1] = 2315.8884029414999
$data[1,2] = 4165.4051748204902
$data[1,3] = 2055.3061353142398
$data[1,4] = 4490.9627733667103
$data[2,0] = 3804.5658589551299
$data[2,1] = 2663.8221382521401
$data[2,2] = 5139.2761438796497
$data[2,3] = 2338.56690583206
$data[2,4] = 5606.76814133065
$data[3,0] = 4295.1131265444601
$data[3,1] = 2928.2062798737902
$data[3,2] = 5972.3332521664197
$data[3,3] = 2549.8265668353902
$data[3,4] = 6579.1010835063998
$data[4,0] = 4711.58796421016
$data[4,1] = 3142.7182742759801
$data[4,2] = 6712.3480138505802
$data[4,3] = 2718.8076697487299
$data[4,4] = 7455.19463625148
$data[5,0] = 5076.3908391831901
$data[5,1] = 3323.8281653497702
$data[5,2] = 7384.8144455297597
$data[5,3] = 2859.8352475247002
$data[5,4] = 8260.55859988472
$data[6,0] = 5402.73684310521
$data[6,1] = 3480.90210338527
$data[6,2] = 8005.30889886356
$data[6,3] = 2980.9605668424701
$data[6,4] = 9010.9410715664908
$data[7,0] = 5699.1578631812899
$data[7,1] = 3619.80283937694
$data[7,2] = 8584.1933636623107
$data[7,3] = 3087.1739900514099
$data[7,4] = 9716.9131617450694
$data[8,0] = 5971.5201541107999
$data[8,1] = 3744.4550483497101
$data[8,2] = 9128.7808418940695
$data[8,3] = 3181.7879915631702
$data[8,4] = 10385.995005999401
$data[9,0] = 6224.0499836282297
$data[9,1] = 3857.6203847967699
$data[9,2] = 9644.4631366506692
$data[9,3] = 3267.1161350267798
$data[9,4] = 11023.770595337801
$data[10,0] = 6459.9041166507004
$data[10,1] = 3961.3186622092799
$data[10,2] = 10135.352017072701
$data[10,3] = 3344.8392925176299
$data[10,4] = 11634.5255837023
$data[11,0] = 6681.5094356052005
$data[11,1] = 4057.0734135903599
$data[11,2] = 10604.6685704667
$data[11,3] = 3416.2177421309998
$data[11,4] = 12221.6368687757
$data[12,0] = 6890.7761984919198
$data[12,1] = 4146.0632011356302
$data[12,2] = 11054.992166926
$data[12,3] = 3482.2210740588798
$data[12,4] = 12787.823053619901
$data[13,0] = 7089.2378469999703
$data[13,1] = 4229.2191148000002
$data[13,2] = 11488.426424779
$data[13,3] = 3543.61143940792
$data[13,4] = 13335.3122508694
$data[14,0] = 7278.1459844891097
$data[14,1] = 4307.2899635567301
$data[14,2] = 11906.713714384699
$data[14,3] = 3600.9989336291801
$data[14,4] = 13865.9583983917
$data[15,0] = 7458.5368605390004
$data[15,1] = 4380.8872426302196
$data[15,2] = 12311.3164711666
$data[15,3] = 3654.8796208451699
$data[15,4] = 14381.3242214273
$data[16,0] = 7631.2791151916099
$data[16,1] = 4450.5169860476799
$data[16,2] = 12703.4763780907
$data[16,3] = 3705.6623518072101
$data[16,4] = 14882.741863950099
$data[17,0] = 7797.1088311145604
$data[17,1] = 4516.6028536846097
$data[17,2] = 13084.2583663121
$data[17,3] = 3753.6881231491702
$data[17,4] = 15371.358141222399
$data[18,0] = 7956.6557688509101
$data[18,1] = 4579.5032040918604
$data[18,2] = 13454.583941081
$data[18,3] = 3799.2443391143102
$data[18,4] = 15848.168939306899
$data[19,0] = 8110.4633394623897
$data[19,1] = 4639.5239450135496
$data[19,2] = 13815.256838388699
$data[19,3] = 3842.57550774011
$data[19,4] = 16314.045790005101
$data[20,0] = 8259.0040406778207
$data[20,1] = 4696.9283588928802
$data[20,2] = 14166.9830657781
$data[20,3] = 3883.8913914581299
$data[20,4] = 16769.7566972282
$data[21,0] = 8402.6915489842304
$data[21,1] = 4751.9447216820599
$data[21,2] = 14510.386760739
$data[21,3] = 3923.3733068565698
$data[21,4] = 17215.982668585701
$data[22,0] = 8541.8903077712494
$data[22,1] = 4804.7722856585997
$data[22,2] = 14846.0228866571
$data[22,3] = 3961.17905657335
$data[22,4] = 17653.330989774699
$data[23,0] = 8676.9232139456799
$data[23,1] = 4855.5860315436103
$data[23,2] = 15174.387504656101
$data[23,3] = 3997.4468352601598
$data[23,4] = 18082.345995045602
$data[24,0] = 8808.0778419242106
$data[24,1] = 4904.5404825002897
$data[24,2] = 15495.9261641665
$data[24,3] = 4032.2983557419798
$data[24,4] = 18503.517889088798
$data[25,0] = 8935.6115294361207
$data[25,1] = 4951.7727943827404
$data[25,2] = 15811.0408169605
$data[25,3] = 4065.8413752032302
$data[25,4] = 18917.2900355196
$data[26,0] = 9059.7555681306803
$data[26,1] = 4997.4052814446904
$data[26,2] = 16120.095560329601
$data[26,3] = 4098.1717546030204
$data[26,4] = 19324.065026318302
$data[27,0] = 9180.7186831959407
$data[27,1] = 5041.54749722625
$data[27,2] = 16423.421442983901
$data[27,3] = 4129.3751512256204
$data[27,4] = 19724.2097730356
$data[28,0] = 9298.6899431873408
$data[28,1] = 5084.2979616736902
$data[28,2] = 16721.320514100102
$data[28,3] = 4159.5284201655204
$data[28,4] = 20118.059806210302
$data[29,0] = 9413.8412094102805
$data[29,1] = 5125.7456044783903
$data[29,2] = 17014.0692562767
$data[29,3] = 4188.7007828706901
$data[29,4] = 20505.922928794698
$data[30,0] = 9526.3292103368294
$data[30,1] = 5165.9709789536801
$data[30,2] = 17301.9215132235
$data[30,3] = 4216.9548077522004
$data[30,4] = 20888.082338629101
$data[31,0] = 9636.2973084714504
$data[31,1] = 5205.0472889906196
$data[31,2] = 17585.1110001979
$data[31,3] = 4244.3472380329003
$data[31,4] = 21264.7993115227
$data[32,0] = 9743.8770132683294
$data[32,1] = 5243.0412626908001
$data[32,2] = 17863.853467631601
$data[32,3] = 4270.9296945543501
$data[32,4] = 21636.315518373602
$data[33,0] = 9849.1892830474899
$data[33,1] = 5280.0138994208201
$data[33,2] = 18138.3485747591
$data[33,3] = 4296.7492755619596
$data[33,4] = 22002.8550356667
$data[34,0] = 9952.3456505654194
$data[34,1] = 5316.0211117335502
$data[34,2] = 18408.781519374199
$data[34,3] = 4321.8490710898104
$data[34,4] = 22364.626097623499
$data[35,0] = 10053.449200392801
$data[35,1] = 5351.1142794714497
$data[35,2] = 18675.324461416101
$data[35,3] = 4346.2686061456398
$data[35,4] = 22721.822629535101
$data[36,0] = 10152.595421113099
$data[36,1] = 5385.3407301223297
$data[36,2] = 18938.137771385798
$data[36,3] = 4370.0442242139998
$data[36,4] = 23074.625594841
$data[37,0] = 10249.8729512667
$data[37,1] = 5418.7441569316297
$data[37,2] = 19197.371129231698
$data[37,3] = 4393.2094204774103
$data[37,4] = 23423.204182936301
$data[38,0] = 10345.364234689099
$data[38,1] = 5451.3649842306904
$data[38,2] = 19453.1644950228
$data[38,3] = 4415.7951324713304
$data[38,4] = 23767.716860177101
$data[39,0] = 10439.1460982505
$data[39,1] = 5483.2406878013999
$data[39,2] = 19705.648969226499
$data[39,3] = 4437.8299945408698
$data[39,4] = 24108.3123029009
$data[40,0] = 10531.290262865999
$data[40,1] = 5514.4060767762103
$data[40,2] = 19954.947557557
$data[40,3] = 4459.3405613821496
$data[40,4] = 24445.1302282865
$data[41,0] = 10621.8637968943
$data[41,1] = 5544.89354250006
$data[41,2] = 20201.175853016
$data[41,3] = 4480.35150507243
$data[41,4] = 24778.302136427999
$data[42,0] = 10710.9295196188
$data[42,1] = 5574.7332789064203
$data[42,2] = 20444.4426458243
$data[42,3] = 4500.8857892778497
$data[42,4] = 25107.9519749723
$data[43,0] = 10798.546361319
$data[43,1] = 5603.9534782431401
$data[43,2] = 20684.8504703459
$data[43,3] = 4520.9648237417296
$data[43,4] = 25434.1967359917
$data[44,0] = 10884.769685470699
$data[44,1] = 5632.5805053931999
$data[44,2] = 20922.4960967794
$data[44,3] = 4540.6086016750896
$data[44,4] = 25757.146993366201
$data[45,0] = 10969.6515777965
$data[45,1] = 5660.6390535466599
$data[45,2] = 21157.470974281601
$data[45,3] = 4559.8358222727302
$data[45,4] = 26076.907387778301
$data[46,0] = 11053.2411062145
$data[46,1] = 5688.1522835742499
$data[46,2] = 21389.861631260599
$data[46,3] = 4578.6640002479999
$data[46,4] = 26393.577065443202
$data[47,0] = 11135.584555162801
$data[47,1] = 5715.1419491137203
$data[47,2] = 21619.750037788701
$data[47,3] = 4597.1095640039302
$data[47,4] = 26707.250075866301
$data[48,0] = 11216.725637298699
$data[48,1] = 5741.6285090957399
$data[48,2] = 21847.213934425701
$data[48,3] = 4615.1879438277601
$data[48,4] = 27018.015733218399
$data[49,0] = 11296.705685168699
$data[49,1] = 5767.6312291971999
$data[49,2] = 22072.3271311797
$data[49,3] = 4632.9136513023795
$data[49,4] = 27325.9589453238
$data[50,0] = 11375.563825101201
$data[50,1] = 5793.1682735076502
$data[50,2] = 22295.159779852998
$data[50,3] = 4650.3003509646096
$data[50,4] = 27631.1605137447
$data[51,0] = 11453.337135281599
$data[51,1] = 5818.2567875232698
$data[51,2] = 22515.778622611699
$data[51,3] = 4667.3609251020198
$data[51,4] = 27933.697408014399
$data[52,0] = 11530.0607897222
$data[52,1] = 5842.9129734375301
$data[52,2] = 22734.247219267101
$data[52,3] = 4684.1075324625099
$data[52,4] = 28233.643016692698
$data[53,0] = 11605.768189623401
$data[53,1] = 5867.1521585730097
$data[53,2] = 22950.626155456001
$data[53,3] = 4700.5516615508805
$data[53,4] = 28531.067377600601
$data[54,0] = 11680.4910834413
$data[54,1] = 5890.9888576929097
$data[54,2] = 23164.973233642399
$data[54,3] = 4716.7041791008196
$data[54,4] = 28826.0373893094
$data[55,0] = 11754.2596768165
$data[55,1] = 5914.4368298393401
$data[55,2] = 23377.343648644499
$data[55,3] = 4732.57537423756
$data[55,4] = 29118.6170057196
$data[56,0] = 11827.1027333853
$data[56,1] = 5937.5091302667197
$data[56,2] = 23587.790149189401
$data[56,3] = 4748.17499878325
$data[56,4] = 29408.867415357301
$data[57,0] = 11899.047667373001
$data[57,1] = 5960.2181579712396
$data[57,2] = 23796.363186832001
$data[57,3] = 4763.51230410277
$data[57,4] = 29696.847206833001
$data[58,0] = 11970.1206287679
$data[58,1] = 5982.57569925799
$data[58,2] = 24003.111053424502
$data[58,3] = 4778.5960748405196
$data[58,4] = 29982.6125217472
$data[59,0] = 12040.346581785299
$data[59,1] = 6004.5929677369504
$data[59,2] = 24208.0800081938
$data[59,3] = 4793.4346598581697
$data[59,4] = 30266.217196192101
$data[60,0] = 12109.7493772503
$data[60,1] = 6026.2806410942203
$data[60,2] = 24411.314395370799
$data[60,3] = 4808.0360006474803
$data[60,4] = 30547.712891870498
$data[61,0] = 12178.3518194637
$data[61,1] = 6047.6488949464601
$data[61,2] = 24612.8567532163
$data[61,3] = 4822.4076574619603
$data[61,4] = 30827.149217753798
$data[62,0] = 12246.175728050401
$data[62,1] = 6068.7074340526196
$data[62,2] = 24812.747915199001
$data[62,3] = 4836.55683338364
$data[62,4] = 31104.5738430987
$data[63,0] = 12313.241995241
$data[63,1] = 6089.4655211274303
$data[63,2] = 25011.027104006302
$data[63,3] = 4850.4903965180802
$data[63,4] = 31380.032602564501
$data[64,0] = 12379.5706389886
$data[64,1] = 6109.9320034748098
$data[64,2] = 25207.732018997001
$data[64,3] = 4864.21490048943
$data[64,4] = 31653.569594095701
$data[65,0] = 12445.180852282299
$data[65,1] = 6130.1153376370003
$data[65,2] = 25402.898917647399
$data[65,3] = 4877.7366033900398
$data[65,4] = 31925.227270171399
$data[66,0] = 12510.0910489831
$data[66,1] = 6150.0236122345796
$data[66,2] = 25596.562691487401
$data[66,3] = 4891.0614853219904
$data[66,4] = 32195.046522962701
$data[67,0] = 12574.3189064736
$data[67,1] = 6169.6645691547401
$data[67,2] = 25788.756936973401
$data[67,3] = 4904.1952646547197
$data[67,4] = 32463.066763889499
$data[68,0] = 12637.8814053884
$data[68,1] = 6189.0456232297201
$data[68,2] = 25979.5140217074
$data[68,3] = 4917.1434131098204
$data[68,4] = 32729.325998021999
$data[69,0] = 12700.7948666623
$data[69,1] = 6208.1738805328496
$data[69,2] = 26168.865146367702
$data[69,3] = 4929.9111697732396
$data[69,4] = 32993.860893728299
$data[70,0] = 12763.0749861124
$data[70,1] = 6227.0561554075402
$data[70,2] = 26356.840402686801
$data[70,3] = 4942.50355412513
$data[70,4] = 33256.706847937203
$data[71,0] = 12824.736866752301
$data[71,1] = 6245.6989863333301
$data[71,2] = 26543.468827780001
$data[71,3] = 4954.9253781689004
$data[71,4] = 33517.8980473466
$data[72,0] = 12885.795049013999
$data[72,1] = 6264.1086507233504
$data[72,2] = 26728.778455102802
$data[72,3] = 4967.1812577334804
$data[72,4] = 33777.467525885797
$data[73,0] = 12946.2635390411
$data[73,1] = 6282.2911787384801
$data[73,2] = 26912.796362286001
$data[73,3] = 4979.2756230150999
$data[73,4] = 34035.447218703899
$data[74,0] = 13006.155835199501
$data[74,1] = 6300.2523661960204
$data[74,2] = 27095.5487160837
$data[74,3] = 4991.2127284197304
$data[74,4] = 34291.868012942497
$data[75,0] = 13065.4849529407
$data[75,1] = 6317.99778664323
$data[75,2] = 27277.060814639601
$data[75,3] = 5002.9966617609098
$data[75,4] = 34546.759795521801
$data[76,0] = 13124.2634481395
$data[76,1] = 6335.53280265989
$data[76,2] = 27457.3571272679
$data[76,3] = 5014.6313528629098
$data[76,4] = 34800.151498151601
$data[77,0] = 13182.503439018399
$data[77,1] = 6352.8625764485296
$data[77,2] = 27636.461331922699
$data[77,3] = 5026.1205816151096
$data[77,4] = 35052.071139765001
$data[78,0] = 13240.216626762
$data[78,1] = 6369.9920797655705
$data[78,2] = 27814.396350520499
$data[78,3] = 5037.4679855188597
$data[78,4] = 35302.545866551103
$data[79,0] = 13297.414314912699
$data[79,1] = 6386.9261032421
$data[79,2] = 27991.184382260901
$data[79,3] = 5048.6770667648998
$data[79,4] = 35551.601989750801
$data[80,0] = 13354.1074276375
$data[80,1] = 6403.6692651390504
$data[80,2] = 28166.846935085599
$data[80,3] = 5059.7511988758097
$data[80,4] = 35799.265021368599
$data[81,0] = 13410.3065269408
$data[81,1] = 6420.2260195772797
$data[81,2] = 28341.404855398599
$data[81,3] = 5070.6936329454202
$data[81,4] = 36045.5597079371
$data[82,0] = 13466.0218288991
$data[82,1] = 6436.6006642803704
$data[82,2] = 28514.8783561657
$data[82,3] = 5081.50750350405
$data[82,4] = 36290.510062465502
$data[83,0] = 13521.2632189821
$data[83,1] = 6452.7973478641297
$data[83,2] = 28687.287043497799
$data[83,3] = 5092.1958340361798
$data[83,4] = 36534.139394687198
$data[84,0] = 13576.0402665233
$data[84,1] = 6468.8200767046301
$data[84,2] = 28858.6499418192
$data[84,3] = 5102.7615421751998
$data[84,4] = 36776.470339718602
$data[85,0] = 13630.362238395601
$data[85,1] = 6484.6727214135499
$data[85,2] = 29028.985517710498
$data[85,3] = 5113.2074445974404
$data[85,4] = 37017.524885229599
$data[86,0] = 13684.2381119445
$data[86,1] = 6500.3590229476604
$data[86,2] = 29198.311702508701
$data[86,3] = 5123.5362616362299
$data[86,4] = 37257.324397216602
$data[87,0] = 13737.676587227401
$data[87,1] = 6515.8825983770002
$data[87,2] = 29366.645913747099
$data[87,3] = 5133.7506216351103
$data[87,4] = 37495.889644470502
$data[88,0] = 13790.686098603899
$data[88,1] = 6531.2469463345697
$data[88,2] = 29534.005075502901
$data[88,3] = 5143.8530650575804
$data[88,4] = 37733.240821814397
$data[89,0] = 13843.274825717301
$data[89,1] = 6546.4554521681202
$data[89,2] = 29700.4056377204
$data[89,3] = 5153.8460483694998
$data[89,4] = 37969.3975721881
$data[90,0] = 13895.450703906999
$data[90,1] = 6561.5113928139299
$data[90,2] = 29865.863594574901
$data[90,3] = 5163.7319477093797
$data[90,4] = 38204.379007650103
$data[91,0] = 13947.221434085999
$data[91,1] = 6576.4179414099899
$data[91,2] = 30030.394501931802
$data[91,3] = 5173.5130623599298
$data[91,4] = 38438.203729359397
$data[92,0] = 13998.594492117199
$data[92,1] = 6591.1781716653604
$data[92,2] = 30194.013493954801
$data[92,3] = 5183.1916180339604
$data[92,4] = 38670.889846596001
$data[93,0] = 14049.577137718499
$data[93,1] = 6605.7950620011297
$data[93,2] = 30356.735298916101
$data[93,3] = 5192.7697699862401
$data[93,4] = 38902.454994880201
$data[94,0] = 14100.176422925701
$data[94,1] = 6620.2714994770604
$data[94,2] = 30518.574254252399
$data[94,3] = 5202.24960596251
$data[94,4] = 39132.916353240202
$data[95,0] = 14150.3992001388
$data[95,1] = 6634.6102835170996
$data[95,2] = 30679.544320909899
$data[95,3] = 5211.63314899532
$data[95,4] = 39362.290660672603
$data[96,0] = 14200.2521297769
$data[96,1] = 6648.8141294462703
$data[96,2] = 30839.6590970221
$data[96,3] = 5220.9223600567902
$data[96,4] = 39590.594231850897
$data[97,0] = 14249.741687564099
$data[97,1] = 6662.8856718500301
$data[97,2] = 30998.931830953999
$data[97,3] = 5230.1191405764002
$data[97,4] = 39817.842972113896
$data[98,0] = 14298.8741714684
$data[98,1] = 6676.8274677668096
$data[98,2] = 31157.3754337504
$data[98,3] = 5239.2253348323102
$data[98,4] = 40044.052391778998

# Write the new values into the worksheet in one shot
$ws.Range("A2:E100").Value = $data

# Move the active selection to E1, matching the saved workbook state
$ws.Range("E1").Select()
